$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.342322
$ws.Range("H2").Value = 1.026966
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3659943333333333
$ws.Range("N2").Value = 1.097983
$ws.Range("O2").Value = 0.006726051721149161
$ws.Range("P2").Value = 0.006726051721149162
$ws.Range("Q2").Value = 0.1252879121753333
$ws.Range("R2").Value = 1.127591209578
$ws.Range("S2").Value = 0.006726051721149161
$ws.Range("T2").Value = 0.006726051721149162

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.342322
$ws.Range("H3").Value = 1.026966
$ws.Range("O3").Value = 0.001678071748088335
$ws.Range("P3").Value = 0.001678071748088335
$ws.Range("Q3").Value = 0.03125787824933333
$ws.Range("R3").Value = 0.281320904244
$ws.Range("S3").Value = 0.001678071748088335
$ws.Range("T3").Value = 0.001678071748088335

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.342322
$ws.Range("H4").Value = 1.026966
$ws.Range("M4").Value = 53.897087
$ws.Range("N4").Value = 161.691261
$ws.Range("O4").Value = 0.9904923704135933
$ws.Range("P4").Value = 0.9904923704135934
$ws.Range("Q4").Value = 18.450158616014
$ws.Range("R4").Value = 166.051427544126
$ws.Range("S4").Value = 0.9904923704135933
$ws.Range("T4").Value = 0.9904923704135934

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.342322
$ws.Range("H5").Value = 1.026966
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06004666666666667
$ws.Range("N5").Value = 0.18014
$ws.Range("O5").Value = 0.001103506117169219
$ws.Range("P5").Value = 0.001103506117169219
$ws.Range("Q5").Value = 0.02055529502666667
$ws.Range("R5").Value = 0.18499765524
$ws.Range("S5").Value = 0.001103506117169219
$ws.Range("T5").Value = 0.001103506117169219
